$wb = $excel.ActiveWorkbook

# --- Sheet "Dades_Període": update extraction timestamps, query hour and source URL ---
$ws1 = $wb.Worksheets.Item("Dades_Període")

# Row 2
$ws1.Range("H2").Value = "2026-02-20 09:35:43"
$ws1.Range("I2").Value = "09:00"
$ws1.Range("J2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T09:00Z"

# Rows 3-6 (only DATA_EXTRACCIO changes)
$ws1.Range("H3").Value = "2026-02-20 09:35:45"
$ws1.Range("H4").Value = "2026-02-20 09:35:45"
$ws1.Range("H5").Value = "2026-02-20 09:35:45"
$ws1.Range("H6").Value = "2026-02-20 09:35:45"

# --- Sheet "Estudi_Capçaleres": update source URL ---
$ws2 = $wb.Worksheets.Item("Estudi_Capçaleres")
$ws2.Range("F2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T09:00Z"
